$d = $word.ActiveDocument

# Subject comment cells and their grade cells live in 4 tables (Torah, History, Math, English).
# Each table's row1/col2 is the "vMerge restart" comment cell (1 empty paragraph),
# and row2/col1 holds the grade label paragraph ("ציון:") followed by an empty
# paragraph that carries the numeric grade.
#
# Paragraph indices (discovered via $d.Paragraphs) for the 4 subjects:
#   Torah:   comment = 5,  grade = 8
#   History: comment = 15, grade = 18
#   Math:    comment = 25, grade = 28
#   English: comment = 35, grade = 38

$torahComment = "במחיצת זאת למדנו על חומש ""ויקרא"", למדנו את ההלכות והאיסורים לעומק,`nהייתה אוירת לימוד מצוינת.`nשפרה את תלמידה מקסימה, הרבה בהצלחה!"
$historyComment = "במחצית זאת למדנו על היסטורית השואה, לכל אחת היתה משימה לעשות פרוייקט על השואה, כך שחפרנו עמוק בשורשים.`nשפרה את ילדה נפלאה, הרבה הצלחה!"
$mathComment = "במחצית זאת למדנו על תורת המיספרים הגדולה, התמקדנו על שברים, על תורת המעגל, רדיוס וקטרים, הרחבנו בנושא המשוואות בנעלם אחד,`nשיננו לעצמינו את הכללים החשובים שנזכור לעתיד.`nשפרה את ילדה מצוינת, את מעולה שיהיה לך הרבה הצלחה להמשך!"
$englishComment = "במחצית זאת למדנו את נושא השיכות, והתמקדנו על הדקדוק והזמנים, הווה מושלם והווה פשוט,היתה אוירה טובה!`nנעמי את מצוינת, בהצלחה!`n"

$d.Paragraphs.Item(5).Range.Text = $torahComment
$d.Paragraphs.Item(8).Range.Text = "95"

$d.Paragraphs.Item(15).Range.Text = $historyComment
$d.Paragraphs.Item(18).Range.Text = "92"

$d.Paragraphs.Item(25).Range.Text = $mathComment
$d.Paragraphs.Item(28).Range.Text = "94"

$d.Paragraphs.Item(35).Range.Text = $englishComment
$d.Paragraphs.Item(38).Range.Text = "84"
